$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values below are stored as text (matching source inlineStr cells),
# so force the Text number format before assignment and use a leading
# apostrophe so Excel does not reinterpret numeric/percentage-looking
# strings as numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "'314.92"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "'2.44%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "'39.45"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "'2.23%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "'5.145"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "'0.82%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "'0.08174"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "'0.74%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "'1.970"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "'0.63%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "'8.191"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "'2.96%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "'0.9263"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "'-0.45%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "'0.1404"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "'-0.35%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "'0.1968"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "'0.63%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "'0.09027"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "'-0.22%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "'0.03503"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "'0.09830"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "'0.03%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "'0.001402"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "'-0.28%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "'0.006029"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "'-1.98%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "'3.655"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "'-1.97%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "'4.244"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "'1.10%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "'-5.45%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "'0.3456"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "'-0.21%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "'0.1344"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "'0.20%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "'4.767"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "'-0.64%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "'0.2423"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "'-1.33%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "'0.04367"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "'-1.10%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "'0.001225"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "'0.29%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "'0.004781"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "'-1.11%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "'0.0001301"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "'-0.06%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "'0.0004001"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "'-10.04%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "'0.02156"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "'3.83%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "'0.05200"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "'1.36%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "'0.007617"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "'1.87%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "'0.009843"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "'-3.04%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "'0.1375"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "'1.42%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "'0.002117"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "'-0.74%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "'0.009124"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "'-1.60%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "'0.00006406"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "'2.45%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "'-0.17%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "'0.002766"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "'-8.91%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "'0.001200"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "'-25.08%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "'-0.17%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "'-0.17%"
